# Increment the "Förändrad" (Changed) date in column C by one day
# for every data row on the active worksheet (row 1 is the header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End($xlUp).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)   # Column C ("Förändrad")
    $cell.Value2 = $cell.Value2 + 1
}
